$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 43 (shifts existing rows 43:61 down to 44:62,
# inheriting formatting from the row being pushed down, same as Excel's
# native "Insert Table Rows Above" behaviour).
$ws.Rows("43:43").Insert()

# Populate the newly inserted row with the new metric.
$ws.Range("A43").Value = "Vaccinations"
$ws.Range("B43").Value = 40
$ws.Range("C43").Value = "# New Vaccination Doses (7-day avg)"
$ws.Range("D43").Value = 420
$ws.Range("E43").Value = "X"
$ws.Range("F43").Value = "X"

# The "Metric - Sort" column (D) is a simple incrementing sequence (10 per
# row); fix up the rows that were pushed down so the sequence stays
# consistent (420, 430, 440, ...).
for ($r = 44; $r -le 62; $r++) {
    $ws.Range("D$r").Value = 10 * $r - 10
}

# Grow the table (ListObject) so it covers the newly inserted row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F62"))

# Match the author's final selection location.
$ws.Range("D43").Select()
